# Auto-generated update of Leve profit cached values across all 8 sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 15751
$ws.Range("I74").Value = 12000
$ws.Range("K74").Value = 12000
$ws.Range("M74").Value = -11064
$ws.Range("H77").Value = 15751
$ws.Range("I77").Value = 12000
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55320
$ws.Range("H107").Value = 522
$ws.Range("I107").Value = 423.45456
$ws.Range("K107").Value = 423.45456
$ws.Range("M107").Value = 1496.54544
$ws.Range("H116").Value = 5420.3335
$ws.Range("I116").Value = 3891.2
$ws.Range("K116").Value = 3891.2
$ws.Range("M116").Value = -449.1999999999998
$ws.Range("H132").Value = 2577.7646
$ws.Range("I132").Value = 2363.875
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 7091.625
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -4561.625
$ws.Range("N132").Value = -23060
$ws.Range("H138").Value = 3606.4707
$ws.Range("J138").Value = 4758.4546
$ws.Range("L138").Value = 14275.3638
$ws.Range("N138").Value = -24555.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4433.899
$ws.Range("I32").Value = 4042.0513
$ws.Range("J32").Value = 34998
$ws.Range("K32").Value = 4042.0513
$ws.Range("L32").Value = 34998
$ws.Range("M32").Value = -3755.0513
$ws.Range("N32").Value = -35572
$ws.Range("H36").Value = 13902.5
$ws.Range("I36").Value = 6006.25
$ws.Range("J36").Value = 19166.666
$ws.Range("K36").Value = 6006.25
$ws.Range("L36").Value = 19166.666
$ws.Range("M36").Value = -5660.25
$ws.Range("N36").Value = -19858.666
$ws.Range("H45").Value = 2226.5
$ws.Range("I45").Value = 2226.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 2226.5
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1849.5
$ws.Range("N45").ClearContents()
$ws.Range("H61").Value = 5760.143
$ws.Range("I61").Value = 4662.3335
$ws.Range("J61").Value = 8504.666999999999
$ws.Range("K61").Value = 4662.3335
$ws.Range("L61").Value = 8504.666999999999
$ws.Range("M61").Value = -4450.3335
$ws.Range("N61").Value = -8928.666999999999
$ws.Range("H110").Value = 1715.9395
$ws.Range("I110").Value = 1208.2069
$ws.Range("K110").Value = 1208.2069
$ws.Range("M110").Value = 836.7931000000001
$ws.Range("H122").Value = 2008.3115
$ws.Range("I122").Value = 1652.878
$ws.Range("J122").Value = 2736.95
$ws.Range("K122").Value = 4958.634
$ws.Range("L122").Value = 8210.849999999999
$ws.Range("M122").Value = -2508.634
$ws.Range("N122").Value = -13110.85
$ws.Range("H132").Value = 2380.0938
$ws.Range("I132").Value = 913.96
$ws.Range("K132").Value = 2741.88
$ws.Range("M132").Value = -211.8800000000001
$ws.Range("H136").Value = 5760.143
$ws.Range("I136").Value = 4662.3335
$ws.Range("J136").Value = 8504.666999999999
$ws.Range("K136").Value = 13987.0005
$ws.Range("L136").Value = 25514.001
$ws.Range("M136").Value = -11437.0005
$ws.Range("N136").Value = -30614.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5835.1577
$ws.Range("I86").Value = 5385.2812
$ws.Range("J86").Value = 8234.5
$ws.Range("K86").Value = 5385.2812
$ws.Range("L86").Value = 8234.5
$ws.Range("M86").Value = -4262.2812
$ws.Range("N86").Value = -10480.5
$ws.Range("H89").Value = 5835.1577
$ws.Range("I89").Value = 5385.2812
$ws.Range("J89").Value = 8234.5
$ws.Range("K89").Value = 26926.406
$ws.Range("L89").Value = 41172.5
$ws.Range("M89").Value = -21310.406
$ws.Range("N89").Value = -52404.5
$ws.Range("H99").Value = 2271.1428
$ws.Range("I99").Value = 1983
$ws.Range("K99").Value = 1983
$ws.Range("M99").Value = -485
$ws.Range("H134").Value = 1600.7142
$ws.Range("I134").Value = 1067.3721
$ws.Range("K134").Value = 3202.1163
$ws.Range("M134").Value = -667.1163000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38134.29
$ws.Range("J31").Value = 46646.6
$ws.Range("L31").Value = 46646.6
$ws.Range("N31").Value = -47236.6
$ws.Range("H34").Value = 38134.29
$ws.Range("J34").Value = 46646.6
$ws.Range("L34").Value = 46646.6
$ws.Range("N34").Value = -47050.6
$ws.Range("H58").Value = 3715.724
$ws.Range("J58").Value = 5354.4287
$ws.Range("L58").Value = 5354.4287
$ws.Range("N58").Value = -5760.4287
$ws.Range("H86").Value = 3567.158
$ws.Range("I86").Value = 3280.9412
$ws.Range("K86").Value = 3280.9412
$ws.Range("M86").Value = -2157.9412
$ws.Range("H89").Value = 3567.158
$ws.Range("I89").Value = 3280.9412
$ws.Range("K89").Value = 16404.706
$ws.Range("M89").Value = -10788.706
$ws.Range("H99").Value = 1823.2632
$ws.Range("I99").Value = 1547.1305
$ws.Range("K99").Value = 1547.1305
$ws.Range("M99").Value = -49.13049999999998
$ws.Range("H122").Value = 2437.8865
$ws.Range("I122").Value = 977.6923
$ws.Range("J122").Value = 4547.0557
$ws.Range("K122").Value = 2933.0769
$ws.Range("L122").Value = 13641.1671
$ws.Range("M122").Value = -483.0769
$ws.Range("N122").Value = -18541.1671
$ws.Range("H126").Value = 1823.2632
$ws.Range("I126").Value = 1547.1305
$ws.Range("K126").Value = 4641.3915
$ws.Range("M126").Value = -2171.3915
$ws.Range("H132").Value = 4499.923
$ws.Range("I132").Value = 3899.2942
$ws.Range("J132").Value = 5634.4443
$ws.Range("K132").Value = 11697.8826
$ws.Range("L132").Value = 16903.3329
$ws.Range("M132").Value = -9167.882599999999
$ws.Range("N132").Value = -21963.3329
$ws.Range("H134").Value = 2956.44
$ws.Range("I134").Value = 1605.1177
$ws.Range("J134").Value = 5828
$ws.Range("K134").Value = 4815.3531
$ws.Range("L134").Value = 17484
$ws.Range("M134").Value = -2280.3531
$ws.Range("N134").Value = -22554
$ws.Range("H136").Value = 3715.724
$ws.Range("J136").Value = 5354.4287
$ws.Range("L136").Value = 16063.2861
$ws.Range("N136").Value = -21163.2861
$ws.Range("H141").Value = 127059.84
$ws.Range("J141").Value = 133481.5
$ws.Range("L141").Value = 133481.5
$ws.Range("N141").Value = -143841.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 4000
$ws.Range("I101").Value = 4000
$ws.Range("K101").Value = 12000
$ws.Range("M101").Value = -9566
$ws.Range("H136").Value = 17546798
$ws.Range("I136").Value = 23811120
$ws.Range("K136").Value = 71433360
$ws.Range("M136").Value = -71428260
$ws.Range("H138").Value = 1556.4286
$ws.Range("I138").Value = 1248.75
$ws.Range("J138").Value = 1966.6666
$ws.Range("K138").Value = 3746.25
$ws.Range("L138").Value = 5899.9998
$ws.Range("M138").Value = 1393.75
$ws.Range("N138").Value = -16179.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 34670
$ws.Range("J123").Value = 34670
$ws.Range("L123").Value = 34670
$ws.Range("N123").Value = -39570
$ws.Range("H126").Value = 2298.0698
$ws.Range("I126").Value = 1398.5518
$ws.Range("J126").Value = 4161.357
$ws.Range("K126").Value = 4195.6554
$ws.Range("L126").Value = 12484.071
$ws.Range("M126").Value = -1725.6554
$ws.Range("N126").Value = -17424.071
$ws.Range("H132").Value = 2048.3918
$ws.Range("I132").Value = 1800.5652
$ws.Range("J132").Value = 5468.4
$ws.Range("K132").Value = 5401.6956
$ws.Range("L132").Value = 16405.2
$ws.Range("M132").Value = -2871.6956
$ws.Range("N132").Value = -21465.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 10480.2
$ws.Range("J46").Value = 16668
$ws.Range("L46").Value = 16668
$ws.Range("N46").Value = -17044
$ws.Range("H136").Value = 4577.1875
$ws.Range("I136").Value = 3314.7222
$ws.Range("K136").Value = 9944.1666
$ws.Range("M136").Value = -7394.1666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 18500
$ws.Range("I25").Value = 7000
$ws.Range("J25").Value = 30000
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = -6707
$ws.Range("N25").Value = -30586
$ws.Range("H43").Value = 25000
$ws.Range("I43").Value = 25000
$ws.Range("K43").Value = 25000
$ws.Range("M43").Value = -24851
$ws.Range("H75").Value = 49997.5
$ws.Range("I75").Value = 49997.5
$ws.Range("K75").Value = 49997.5
$ws.Range("M75").Value = -49061.5
$ws.Range("H78").Value = 49997.5
$ws.Range("I78").Value = 49997.5
$ws.Range("K78").Value = 149992.5
$ws.Range("M78").Value = -145312.5
$ws.Range("H103").Value = 42999.668
$ws.Range("J103").Value = 42999.668
$ws.Range("L103").Value = 42999.668
$ws.Range("N103").Value = -45343.668
$ws.Range("H117").Value = 58140
$ws.Range("J117").Value = 58140
$ws.Range("L117").Value = 58140
$ws.Range("N117").Value = -67318
$ws.Range("H126").Value = 2100.0625
$ws.Range("I126").Value = 1548.08
$ws.Range("J126").Value = 4161.357
$ws.Range("K126").Value = 4644.24
$ws.Range("L126").Value = 12484.24
$ws.Range("M126").Value = -2174.24
